$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update computed values (re-run model outputs) ---
$ws.Range("D2").Value = 0.082146520538613807
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 57188727.426597431

$ws.Range("D3").Value = 0.081540471185775315
$ws.Range("E3").Value = 0.000000000000000029695880345295601
$ws.Range("F3").Value = 57220263.45451048

$ws.Range("D4").Value = 0.082146520538470588
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 57188727.426606447

$ws.Range("D5").Value = 0.082146520538367726
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 57188727.426604457

$ws.Range("D6").Value = 0.081540471185775301
$ws.Range("E6").Value = 0.00000000000000003083802958934543
$ws.Range("F6").Value = 57220263.45451048

$ws.Range("D7").Value = 0.081540471185775357
$ws.Range("E7").Value = 0.000000000000000023985134125046449
$ws.Range("F7").Value = 57220263.45451048

$ws.Range("D8").Value = 0.1267959660359094
$ws.Range("E8").Value = 0.00000000000000001941653714884712
$ws.Range("F8").Value = 52961462.572426662

# --- Apply number formats (named cell styles: Currency then Percent, so
#     Currency ends up as xfId 1 / cellXfs 2 and Percent as xfId 2 / cellXfs 3,
#     matching the target workbook's style ordering) ---
$ws.Range("F2:F8").Style = "Currency"
$ws.Range("D2:E8").Style = "Percent"

# --- Column width / sheet view cosmetics ---
$ws.Columns("F").ColumnWidth = 13.8
$null = $ws.Range("B2:F8").Select()
